# adding the new test cases
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Credentials: add a new row of test data (testuser2 / 123)
# ---------------------------------------------------------------------------
$credentials = $wb.Worksheets.Item("Credentials")
$credentials.Range("A4").Value = "testuser2"
$credentials.Range("B4").Value = 123

# ---------------------------------------------------------------------------
# 2. Create the three new worksheets. They are appended after the last
#    existing sheet (CustomActions) and moved into their final position
#    afterwards - this also reproduces the sheetId numbering Excel assigned
#    in the real edit (BaseLine_Creation=7, ChangeNotice=6, Product_Creation=9)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$changeNotice = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$changeNotice.Name = "ChangeNotice"

$baseLineCreation = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$baseLineCreation.Name = "BaseLine_Creation"

$placeholder = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$placeholder.Name = "zzTempPlaceholder"

$productCreation = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$productCreation.Name = "Product_Creation"

$toRemove = $wb.Worksheets.Item("zzTempPlaceholder")
$toRemove.Delete()

# ---------------------------------------------------------------------------
# 3. BaseLine_Creation content
# ---------------------------------------------------------------------------
$baseLine = $wb.Worksheets.Item("BaseLine_Creation")
$baseLine.Columns.Item(1).ColumnWidth = 14.166666666666666
$baseLine.Columns.Item(2).ColumnWidth = 20.276041666666668

$baseLine.Range("A1").Value = "Name"
$baseLine.Range("B1").Value = "Description"
$baseLine.Range("A2").Value = "BaseLine 1"
$baseLine.Range("B2").Value = "For testing"
$baseLine.Range("A3").Value = "BaseLine 2"
$baseLine.Range("B3").Value = "For testing"

$baseLine.PageSetup.Orientation = 1
$baseLine.Range("A3").Select()

# ---------------------------------------------------------------------------
# 4. ChangeNotice content
# ---------------------------------------------------------------------------
$changeNotice = $wb.Worksheets.Item("ChangeNotice")
$changeNotice.Columns.Item(1).ColumnWidth = 24.721354166666668
$changeNotice.Columns.Item(2).ColumnWidth = 22.166666666666668
$changeNotice.Columns.Item(3).ColumnWidth = 27.385416666666668
$changeNotice.Columns.Item(4).ColumnWidth = 21.498697916666668

$changeNotice.Range("A1").Value = "Name_CN"
$changeNotice.Range("B1").Value = "Name_CT"
$changeNotice.Range("C1").Value = "Approver"
$changeNotice.Range("D1").Value = "Reviewer"
$changeNotice.Range("A2").Value = "ChangeNotice_1"
$changeNotice.Range("B2").Value = "ChangeTask_1"
$changeNotice.Range("C2").Value = "testuser2"
$changeNotice.Range("D2").Value = "testuser2"

$changeNotice.PageSetup.Orientation = 1
$changeNotice.Range("A2").Select()

# ---------------------------------------------------------------------------
# 5. Product_Creation content
# ---------------------------------------------------------------------------
$productCreation = $wb.Worksheets.Item("Product_Creation")
$productCreation.Columns.Item(1).ColumnWidth = 21.944010416666668
$productCreation.Columns.Item(2).ColumnWidth = 17.721354166666668

$productCreation.Range("A1").Value = "Name"
$productCreation.Range("B1").Value = "Description"
$productCreation.Range("A2").Value = "Product1"
$productCreation.Range("B2").Value = "For Testing "
$productCreation.Range("A3").Value = "Product2"
$productCreation.Range("B3").Value = "For Testing "

$productCreation.PageSetup.Orientation = 1
$productCreation.Range("B4").Select()

# ---------------------------------------------------------------------------
# 6. CustomActions: row 2 loses its explicit custom height (back to default)
# ---------------------------------------------------------------------------
$customActions = $wb.Worksheets.Item("CustomActions")
$customActions.Rows.Item(2).AutoFit()

# ---------------------------------------------------------------------------
# 7. Move the new sheets into their final tab order:
#    Credentials, Part_Creation, BaseLine_Creation, ContentFileManagement,
#    ReportsManagement, CustomActions, ChangeNotice, Product_Creation
# ---------------------------------------------------------------------------
$baseLineRef = $wb.Worksheets.Item("BaseLine_Creation")
$afterPart = $wb.Worksheets.Item("Part_Creation")
$baseLineRef.Move([Type]::Missing, $afterPart)

$changeNoticeRef = $wb.Worksheets.Item("ChangeNotice")
$afterCustomActions = $wb.Worksheets.Item("CustomActions")
$changeNoticeRef.Move([Type]::Missing, $afterCustomActions)

$productCreationRef = $wb.Worksheets.Item("Product_Creation")
$afterChangeNotice = $wb.Worksheets.Item("ChangeNotice")
$productCreationRef.Move([Type]::Missing, $afterChangeNotice)

# ---------------------------------------------------------------------------
# 8. Credentials selection (selected cell moves to B4, tab no longer active)
# ---------------------------------------------------------------------------
$credentials = $wb.Worksheets.Item("Credentials")
$credentials.Range("B4").Select()

# ---------------------------------------------------------------------------
# 9. Product_Creation becomes the active sheet/tab
# ---------------------------------------------------------------------------
$productCreation = $wb.Worksheets.Item("Product_Creation")
$productCreation.Activate()
$productCreation.Range("B4").Select()
